$d = $word.ActiveDocument

# --- Locate the target paragraph ("Iniciativa: ... situação seguinte:") and the
#     empty bookmark-only paragraph that immediately follows it.
$findRng = $d.Content
$null = $findRng.Find.Execute("Esta regra impede empates. O jogo termina em qualquer uma da situação seguinte:", $false)
$targetPara = $d.Range($findRng.Start, $findRng.Start).Paragraphs(1)
$bookmarkPara = $targetPara.Next()

$replaceRange = $d.Range($targetPara.Range.Start, $bookmarkPara.Range.End)

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rStyle w:val="SubtleEmphasis"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Iniciativa</w:t></w:r><w:r w:rsidRPr="001E507B"><w:rPr><w:rStyle w:val="SubtleEmphasis"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:rStyle w:val="SubtleEmphasis"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Esta regra impede empates. O jogo termina em qualquer uma </w:t></w:r><w:r><w:t>das seguintes situações</w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/><w:r><w:t>:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Ambos os jogadores completam uma ilha ao mesmo tempo.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Os jogadores passam em turnos consecutivos.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Um jogador passa 4 vezes consecutivas.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>Nesta situação, o jogador vencedor é aquele que tem mais peças afundadas desde que o outro jogador tenha afundado alguma peça. Numa situação improvável em que nenhuma peça tenha sido afundada, ganha o jogador que jogou primeiro.</w:t></w:r></w:p>
'@
$replaceRange.InsertXML($newXml)

# --- Turn the three new plain paragraphs into a proper bulleted list.
$p1Rng = $d.Content
$null = $p1Rng.Find.Execute("Ambos os jogadores completam uma ilha ao mesmo tempo.", $false)
$firstItem = $d.Range($p1Rng.Start, $p1Rng.Start).Paragraphs(1)

$p3Rng = $d.Content
$null = $p3Rng.Find.Execute("Um jogador passa 4 vezes consecutivas.", $false)
$thirdItem = $d.Range($p3Rng.Start, $p3Rng.Start).Paragraphs(1)

$listRange = $d.Range($firstItem.Range.Start, $thirdItem.Range.End)
$listRange.Style = "List Paragraph"
$listRange.ListFormat.ApplyBulletDefault()

Write-Host "Edit applied"
